$d = $word.ActiveDocument

# Locate the paragraph that contains "Ver no Jupiter Salvar em pdf Salvar em docx"
# and delete it together with the empty paragraph right before it and the
# copyright paragraph right after it (three paragraphs total get removed).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $prev = $target.Previous()   # the blank paragraph just before it
    $next = $target.Next()       # the copyright paragraph just after it

    $startRange = $prev.Range.Start
    $endRange = $next.Range.End

    $r = $d.Range($startRange, $endRange)
    $r.Delete()
}

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
